# ADD results from server
# Update computed result values on row 2 of each year sheet (2025, 2030, 2035, 2040, 2045, 2050)

$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 9064.015633300272
$ws.Range("E2").Value = 23114.1806051135
$ws.Range("G2").Value = 12143.88856899275
$ws.Range("I2").Value = 51856.62951455999
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 8162.670248152207
$ws.Range("O2").Value = 12103.9598975121

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 6794.473007685274
$ws.Range("B2").Value = 15129.75820914074
$ws.Range("E2").Value = 64447.80433120584
$ws.Range("I2").Value = 96650.97486679791
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 13869.23146919471
$ws.Range("O2").Value = 11702.01983856686

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 8850.598356357281
$ws.Range("B2").Value = 23430.21916234464
$ws.Range("E2").Value = 96928.5825258494
$ws.Range("I2").Value = 143329.1971486044
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 25217.17353999086
$ws.Range("O2").Value = 19174.53054577044

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("O2").Value = 1235.461171605062

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 3879.890732167435
$ws.Range("O2").Value = 1020.236817535124

# Sheet "2050" - no changes
